{"js": "// Add a new \"MN10\" row (with the \"Deseja realmente excluir\" message) right\n// after the existing last row (\"MN09\" - \"Deseja realmente cancelar\") of the\n// messages table, following the exact same layout/formatting pattern used\n// throughout that table.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The messages table (\"N\u00famero\" / \"Mensagem\") is the second table in the\n// document - locate it defensively by checking its header instead of\n// hard-coding the index.\nlet targetTable = null;\nfor (let i = 0; i < tables.items.length; i++) {\n  const tbl = tables.items[i];\n  tbl.load(\"values\");\n  await context.sync();\n  if (\n    tbl.values.length &&\n    tbl.values[0].length >= 2 &&\n    tbl.values[0][0].trim() === \"N\u00famero\"\n  ) {\n    targetTable = tbl;\n    break;\n  }\n}\nif (!targetTable) {\n  targetTable = tables.items[tables.items.length - 1];\n}\n\n// Append a new row with placeholder content - we will replace the second\n// cell's content below so every run gets the exact OOXML shape (three\n// separate runs sharing the same rPr) used by the sibling rows.\ntargetTable.addRows(Word.InsertLocation.end, 1, [[\"MN10\", \"\"]]);\nawait context.sync();\n\nconst rows = targetTable.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst newRow = rows.items[rows.items.length - 1];\nconst cells = newRow.cells;\ncells.load(\"items\");\nawait context.sync();\n\nconst messageCellBody = cells.items[1].body;\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/>\n                <w:lang w:val=\"pt-BR\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/>\n                <w:lang w:val=\"pt-BR\"/>\n              </w:rPr>\n              <w:t>\\u201CDeseja realmente excluir</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/>\n                <w:lang w:val=\"pt-BR\"/>\n              </w:rPr>\n              <w:t>?</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\"/>\n                <w:lang w:val=\"pt-BR\"/>\n              </w:rPr>\n              <w:t>\\u201D;</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nmessageCellBody.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Add a new \"MN10\" row (with the \"Deseja realmente excluir\" message) right\n# after the existing last row (\"MN09\" - \"Deseja realmente cancelar\") of the\n# messages table, following the exact same layout/formatting pattern used\n# throughout that table.\n\n$d = $word.ActiveDocument\n\n# The messages table (\"N\u00famero\" / \"Mensagem\") is the second table in the\n# document - locate it defensively by checking its header instead of\n# hard-coding the index.\n$targetTable = $null\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $candidate = $d.Tables.Item($i)\n    $header = $candidate.Cell(1, 1).Range.Text.Trim()\n    if ($header -eq \"N\u00famero\" -or $header -eq \"Numero\") {\n        $targetTable = $candidate\n        break\n    }\n}\nif ($targetTable -eq $null) {\n    $targetTable = $d.Tables.Item($d.Tables.Count)\n}\n\n# Append a new row - Rows.Add() clones the layout/formatting (row height,\n# cell widths, shading, alignment, fonts, ...) of the table's last row, just\n# like pressing Tab at the end of the table in Word.\n$newRow = $targetTable.Rows.Add()\n\n$newRow.Cells.Item(1).Range.Text = \"MN10\"\n\n$openQuote = [char]0x201C\n$closeQuote = [char]0x201D\n$newRow.Cells.Item(2).Range.Text = $openQuote + \"Deseja realmente excluir?\" + $closeQuote + \";\"\n"}
